# Apply the "Office Theme" design colours to the presentation's theme
# (this is what happens when a different Design is applied from the
# PowerPoint ribbon: the colour scheme bound to the slide master /
# presentation theme - ppt/theme/theme1.xml - is replaced).
#
# Colour order inside a DrawingML <a:clrScheme> (and therefore the
# index used by ThemeColorScheme.Item(n).RGB) is:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#   8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
#
# Target ("Office Theme") values, expressed as the decimal RGB()
# encoding PowerPoint uses internally (R + G*256 + B*65536):
#   dk1      000000 -> 0
#   lt1      FFFFFF -> 16777215
#   dk2      44546A -> 6968388
#   lt2      E7E6E6 -> 15132391
#   accent1  5B9BD5 -> 13998939
#   accent2  ED7D31 -> 3243501
#   accent3  A5A5A5 -> 10855845
#   accent4  FFC000 -> 49407
#   accent5  4472C4 -> 12874308
#   accent6  70AD47 -> 4697456
#   hlink    0563C1 -> 12673797
#   folHlink 954F72 -> 7491477

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme

$officeThemeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $theme.ThemeColorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
